# Daily attendance processing - 2026-01-28 06:50:32
# Normalizes the "Recorded By" (column G) values so that the leading
# "System, " entry is moved to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G whose "Recorded By" value currently starts with "System, "
$rows = @(2,3,4,5,6,7,8,28,29,30,31,32,33,34,54,55,56,57,58,59,60,80,81,82,106,107,108,132,133,134)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value.ToString().StartsWith("System, ")) {
        $rest = $value.ToString().Substring(8)   # strip leading "System, "
        $newValue = $rest + ", System"
        $cell.Value2 = $newValue
    }
}
